# The deck currently has two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (used only by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     colour scheme (used by the Slide Master / all slides)
#
# The authored edit swaps the two themes' content: the design/theme applied
# to the slides becomes the standard "Office" colour scheme (what used to
# live in theme1.xml), while the "Integral" colours move to the other theme
# part. Re-applying the built-in "Office" colour set to the presentation's
# current design reproduces that swap for the theme that actually drives
# the slides' appearance.

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Theme colours (RRGGBB) for the standard Office colour scheme, in the
# PowerPoint ThemeColorScheme index order:
#  1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#  5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink
$officeColors = @(
    (RGB 0x00 0x00 0x00), # 1  dk1
    (RGB 0xFF 0xFF 0xFF), # 2  lt1
    (RGB 0x44 0x54 0x6A), # 3  dk2
    (RGB 0xE7 0xE6 0xE6), # 4  lt2
    (RGB 0x5B 0x9B 0xD5), # 5  accent1
    (RGB 0xED 0x7D 0x31), # 6  accent2
    (RGB 0xA5 0xA5 0xA5), # 7  accent3
    (RGB 0xFF 0xC0 0x00), # 8  accent4
    (RGB 0x44 0x72 0xC4), # 9  accent5
    (RGB 0x70 0xAD 0x47), # 10 accent6
    (RGB 0x05 0x63 0xC1), # 11 hlink
    (RGB 0x95 0x4F 0x72)  # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}

try { $colorScheme.Name = "Office" } catch {}
try { $theme.Name = "Office Theme" } catch {}
